$d = $word.ActiveDocument

$replacements = @(
    @("842×4=", "237×5="),
    @("290×7=", "538×4="),
    @("624×6=", "414×7="),
    @("744×9=", "239×8="),
    @("874×8=", "969×2="),
    @("501×7=", "622×8="),
    @("232×8=", "155×4="),
    @("403×8=", "377×8="),
    @("365×3=", "402×7="),
    @("802×6=", "616×3="),
    @("871×9=", "623×8="),
    @("737×8=", "375×5="),
    @("750×7=", "637×9="),
    @("411×4=", "433×5="),
    @("551×6=", "190×3="),
    @("896×3=", "118×6="),
    @("253×7=", "160×6="),
    @("637×8=", "628×9="),
    @("255×7=", "790×4="),
    @("207×7=", "258×9="),
    @("282×6=", "398×2="),
    @("885×6=", "796×8="),
    @("338×4=", "908×9="),
    @("470×2=", "761×5="),
    @("104×2=", "524×3=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
